# Generate Report for Handback
#
# The handback pipeline just delivered the zh-cn and de-de translations for
# "33916f17-483f-4a27-a3b2-d560772d4bbf.md" (and its content-duplicate
# "80e4ec08-09a8-43f6-8588-8a9c9c1ba20f.md"). Flip their status from
# "Ready for handoff" to "Handed back: in sync with en-US" on the Overview
# sheet, and fill in the Latest Target File / Latest Handback File / Latest
# Handback DateTime columns (plus the new hyperlink on the target file) on
# the per-language detail sheets.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn (col E) and de-de (col F) status for rows 3 and 4
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack
$wsOverview.Range("E4").Value = $statusHandedBack
$wsOverview.Range("F4").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C3").Value = $statusHandedBack
$wsZhCn.Range("J3").Value = "33916f17-483f-4a27-a3b2-d560772d4bbf.b45bc6e6e42c09928d77f7d6a79e2dd528aca2f6.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-31 07:48:39"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/349402c63d0501f6f758630625e459808bfc4d55/e2e/33916f17-483f-4a27-a3b2-d560772d4bbf.md", "", "", "33916f17-483f-4a27-a3b2-d560772d4bbf.md") | Out-Null
$wsZhCn.Range("I3").Font.Underline = 2
$wsZhCn.Range("I3").Font.Color = 15570276

$wsZhCn.Range("C4").Value = $statusHandedBack
$wsZhCn.Range("J4").Value = "33916f17-483f-4a27-a3b2-d560772d4bbf.b45bc6e6e42c09928d77f7d6a79e2dd528aca2f6.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-08-31 07:48:39"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/349402c63d0501f6f758630625e459808bfc4d55/e2e/33916f17-483f-4a27-a3b2-d560772d4bbf.md", "", "", "33916f17-483f-4a27-a3b2-d560772d4bbf.md") | Out-Null
$wsZhCn.Range("I4").Font.Underline = 2
$wsZhCn.Range("I4").Font.Color = 15570276

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C3").Value = $statusHandedBack
$wsDeDe.Range("J3").Value = "33916f17-483f-4a27-a3b2-d560772d4bbf.b45bc6e6e42c09928d77f7d6a79e2dd528aca2f6.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-31 07:48:55"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/5d6cd7220e11814f3250f0d386ea6d27e865554d/e2e/33916f17-483f-4a27-a3b2-d560772d4bbf.md", "", "", "33916f17-483f-4a27-a3b2-d560772d4bbf.md") | Out-Null
$wsDeDe.Range("I3").Font.Underline = 2
$wsDeDe.Range("I3").Font.Color = 15570276

$wsDeDe.Range("C4").Value = $statusHandedBack
$wsDeDe.Range("J4").Value = "33916f17-483f-4a27-a3b2-d560772d4bbf.b45bc6e6e42c09928d77f7d6a79e2dd528aca2f6.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-08-31 07:48:55"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/5d6cd7220e11814f3250f0d386ea6d27e865554d/e2e/33916f17-483f-4a27-a3b2-d560772d4bbf.md", "", "", "33916f17-483f-4a27-a3b2-d560772d4bbf.md") | Out-Null
$wsDeDe.Range("I4").Font.Underline = 2
$wsDeDe.Range("I4").Font.Color = 15570276
